{"js": "// Replace the 25 two-digit-by-two-digit multiplication problems in the\n// document's table with their new values, per the commit diff.\nconst replacements = [\n  [\"12\u00d791=1092\", \"25\u00d747=1175\"],\n  [\"47\u00d792=4324\", \"41\u00d747=1927\"],\n  [\"69\u00d797=6693\", \"23\u00d738=874\"],\n  [\"53\u00d720=1060\", \"86\u00d796=8256\"],\n  [\"31\u00d727=837\", \"71\u00d720=1420\"],\n  [\"45\u00d721=945\", \"65\u00d788=5720\"],\n  [\"80\u00d783=6640\", \"89\u00d782=7298\"],\n  [\"51\u00d759=3009\", \"79\u00d774=5846\"],\n  [\"98\u00d760=5880\", \"59\u00d766=3894\"],\n  [\"93\u00d796=8928\", \"31\u00d749=1519\"],\n  [\"64\u00d790=5760\", \"44\u00d746=2024\"],\n  [\"38\u00d762=2356\", \"36\u00d758=2088\"],\n  [\"11\u00d732=352\", \"87\u00d793=8091\"],\n  [\"78\u00d791=7098\", \"56\u00d768=3808\"],\n  [\"71\u00d778=5538\", \"68\u00d712=816\"],\n  [\"11\u00d713=143\", \"31\u00d777=2387\"],\n  [\"34\u00d746=1564\", \"31\u00d763=1953\"],\n  [\"57\u00d712=684\", \"90\u00d735=3150\"],\n  [\"67\u00d733=2211\", \"28\u00d727=756\"],\n  [\"53\u00d779=4187\", \"20\u00d774=1480\"],\n  [\"64\u00d713=832\", \"66\u00d788=5808\"],\n  [\"94\u00d746=4324\", \"52\u00d752=2704\"],\n  [\"96\u00d797=9312\", \"44\u00d717=748\"],\n  [\"23\u00d791=2093\", \"88\u00d732=2816\"],\n  [\"24\u00d795=2280\", \"54\u00d780=4320\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-by-two-digit multiplication problems in the\n# document's table with their new values, per the commit diff.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"12\u00d791=1092\", \"25\u00d747=1175\"),\n    @(\"47\u00d792=4324\", \"41\u00d747=1927\"),\n    @(\"69\u00d797=6693\", \"23\u00d738=874\"),\n    @(\"53\u00d720=1060\", \"86\u00d796=8256\"),\n    @(\"31\u00d727=837\",  \"71\u00d720=1420\"),\n    @(\"45\u00d721=945\",  \"65\u00d788=5720\"),\n    @(\"80\u00d783=6640\", \"89\u00d782=7298\"),\n    @(\"51\u00d759=3009\", \"79\u00d774=5846\"),\n    @(\"98\u00d760=5880\", \"59\u00d766=3894\"),\n    @(\"93\u00d796=8928\", \"31\u00d749=1519\"),\n    @(\"64\u00d790=5760\", \"44\u00d746=2024\"),\n    @(\"38\u00d762=2356\", \"36\u00d758=2088\"),\n    @(\"11\u00d732=352\",  \"87\u00d793=8091\"),\n    @(\"78\u00d791=7098\", \"56\u00d768=3808\"),\n    @(\"71\u00d778=5538\", \"68\u00d712=816\"),\n    @(\"11\u00d713=143\",  \"31\u00d777=2387\"),\n    @(\"34\u00d746=1564\", \"31\u00d763=1953\"),\n    @(\"57\u00d712=684\",  \"90\u00d735=3150\"),\n    @(\"67\u00d733=2211\", \"28\u00d727=756\"),\n    @(\"53\u00d779=4187\", \"20\u00d774=1480\"),\n    @(\"64\u00d713=832\",  \"66\u00d788=5808\"),\n    @(\"94\u00d746=4324\", \"52\u00d752=2704\"),\n    @(\"96\u00d797=9312\", \"44\u00d717=748\"),\n    @(\"23\u00d791=2093\", \"88\u00d732=2816\"),\n    @(\"24\u00d795=2280\", \"54\u00d780=4320\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
